# Weekly fruit/vegetable price update:
# Insert a new price record as row 130 in the "Vega Modelo de Temuco - Achicoria"
# data block, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at sheet row 130 (shifts rows 130:148 down to 131:149)
$ws.Rows.Item(130).Insert()

# Populate the new row 130 with the new weekly record.
$ws.Range("A130").Value = 10
$ws.Range("B130").Value = "Vega Modelo de Temuco"
$ws.Range("C130").Value = "La Araucanía"
$ws.Range("D130").Value = 45180
$ws.Range("E130").Value = 9
$ws.Range("F130").Value = 100112010
$ws.Range("G130").Value = "Achicoria"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 150
$ws.Range("K130").Value = 10000
$ws.Range("L130").Value = 10000
$ws.Range("M130").Value = 10000
$ws.Range("N130").Value = "$/caja 18 unidades"
$ws.Range("O130").Value = "Región Metropolitana"
$ws.Range("P130").Value = 556
$ws.Range("Q130").Value = 18
$ws.Range("R130").Value = "Hortaliza"
